# Updated cryptos list (price/volume refresh + one row swap for rows 22/23).
# Values are written with a leading apostrophe to force Excel to keep them
# as plain text (matching the original inlineStr/text cell type) instead of
# auto-converting number-like strings (e.g. "5.27") into numeric cells, and
# then the cell style is reset to "Normal" so no stray quote-prefix style
# is left applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '68.425.36'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  +0.44%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '2.648.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  +0.54%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + '  -0.13%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + '597.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  -0.04%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + '159.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  +3.12%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'" + '  +0.00%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'" + '  -0.48%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'" + '  +0.15%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'" + '  -0.80%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'" + '5.27'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  +0.81%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'" + '  +0.73%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + '28.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  +1.39%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'" + '3.132.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  +0.58%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'" + '  -1.40%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '68.353.99'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  +0.54%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + '2.658.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + '  -0.54%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'" + '  +3.15%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'" + '364.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  +0.62%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '7.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  +1.47%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'" + '  +1.97%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'" + 'SuiNetwork'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'" + 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'" + '2.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  +3.16%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'" + 'NEARProtocol'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'" + 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'" + '4.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  -0.70%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'" + '75.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + '  +0.10%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'" + '  +0.09%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + '9.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  +2.21%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E28").Value = "'" + '  -1.80%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'" + '  -0.16%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + '577.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  +2.96%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'" + '  +1.08%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'" + '  +0.91%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'" + '  +1.07%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'" + '  +4.00%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'" + '  +0.08%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  -0.04%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'" + '160.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + '  -0.01%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'" + '19.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  +2.26%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'" + '  -0.33%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'" + '  +0.43%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +0.78%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + '2.64'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  +0.60%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'" + '0.0₆0320'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + '  -5.16%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'" + '  +0.07%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + '158.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  +1.08%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'" + '  +2.50%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + '21.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  +0.95%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'" + '1.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  +1.05%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'" + '  -0.68%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'" + '0.577'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  +3.16%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'" + '0.615'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  +0.14%  '
$ws.Range("E51").Style = "Normal"
